# Update hours of Rifat (column H) in the Tasks worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Add Rifat's hours for the respective tasks (column H).
$ws.Range("H26").Value = 9
$ws.Range("H29").Value = 15
$ws.Range("H30").Value = 12
$ws.Range("H36").Value = 25
$ws.Range("H43").Value = 11
$ws.Range("H51").Value = 8

# Update the view state to match where the edit was made.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("I59").Select()
